$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot (GitHub Actions update).
# Each updated "Price" (column D) cell has NumberFormat forced to text ('@')
# before its value is written, because plain numeric-looking strings (e.g.
# "14.25") would otherwise be auto-converted to real numbers by Excel, while
# the source data stores every Price cell as text (note some prices, like
# "37.397.84", use '.' as a thousands separator and are not valid numbers).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.397.84'
$ws.Range('E2').Value = '  +0.65%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.012.15'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '258.41'
$ws.Range('E5').Value = '  +5.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.617'
$ws.Range('E6').Value = '  -1.44%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '55.95'
$ws.Range('E8').Value = '  -6.40%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.387'
$ws.Range('E9').Value = '  +1.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0770'
$ws.Range('E10').Value = '  -4.24%  '
$ws.Range('E11').Value = '  -1.92%  '
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.310.11'
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.25'
$ws.Range('E13').Value = '  -5.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.804'
$ws.Range('E14').Value = '  -4.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.74'
$ws.Range('E15').Value = '  -7.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.25'
$ws.Range('E16').Value = '  -3.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.004.72'
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.343.16'
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.64'
$ws.Range('E19').Value = '  -0.84%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0840'
$ws.Range('E20').Value = '  -2.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.17'
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '227.60'
$ws.Range('E22').Value = '  -0.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.65'
$ws.Range('E23').Value = '  +7.55%  '
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.34'
$ws.Range('E25').Value = '  -0.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.16'
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.99'
$ws.Range('E27').Value = '  -4.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.63'
$ws.Range('E28').Value = '  +0.16%  '
$ws.Range('E29').Value = '  -8.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.31'
$ws.Range('E30').Value = '  -2.44%  '
$ws.Range('E31').Value = '  -1.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.63'
$ws.Range('E32').Value = '  -3.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0648'
$ws.Range('E33').Value = '  -0.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.53'
$ws.Range('E34').Value = '  +1.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.39'
$ws.Range('E35').Value = '  -0.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.83'
$ws.Range('E36').Value = '  +0.94%  '
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.34'
$ws.Range('E38').Value = '  +1.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.18'
$ws.Range('E39').Value = '  -3.26%  '
$ws.Range('E40').Value = '  +3.62%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.21'
$ws.Range('E41').Value = '  +3.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0938'
$ws.Range('E42').Value = '  -4.12%  '
$ws.Range('E43').Value = '  -0.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.387.44'
$ws.Range('E44').Value = '  +1.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '90.04'
$ws.Range('E45').Value = '  -0.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '15.66'
$ws.Range('E46').Value = '  -6.11%  '
$ws.Range('E47').Value = '  -2.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.09'
$ws.Range('E48').Value = '  -4.79%  '
$ws.Range('E49').Value = '  +1.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.202.81'
$ws.Range('E50').Value = '  +0.67%  '
$ws.Range('E51').Value = '  -4.12%  '
